$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 818.36365
$ws.Range("I12").Value = 794.25
$ws.Range("J12").Value = 882.6667
$ws.Range("K12").Value = 794.25
$ws.Range("L12").Value = 882.6667
$ws.Range("M12").Value = -624.25
$ws.Range("N12").Value = -1222.6667
$ws.Range("H43").Value = 4518.8
$ws.Range("I43").Value = 4264
$ws.Range("K43").Value = 4264
$ws.Range("M43").Value = -4195
$ws.Range("H106").Value = 6999.75
$ws.Range("I106").Value = 6999.75
$ws.Range("K106").Value = 6999.75
$ws.Range("M106").Value = -6368.75
$ws.Range("H116").Value = 3995
$ws.Range("J116").Value = 3995
$ws.Range("L116").Value = 3995
$ws.Range("N116").Value = -10879
$ws.Range("H125").Value = 2176.6667
$ws.Range("I125").Value = 2015
$ws.Range("K125").Value = 18135
$ws.Range("M125").Value = -15675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 5015000
$ws.Range("J24").Value = 5015000
$ws.Range("L24").Value = 5015000
$ws.Range("N24").Value = -5015748
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H100").Value = 5015000
$ws.Range("J100").Value = 5015000
$ws.Range("L100").Value = 5015000
$ws.Range("N100").Value = -5017164

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1143.8334
$ws.Range("I36").Value = 1143.8334
$ws.Range("K36").Value = 1143.8334
$ws.Range("M36").Value = -609.8334
$ws.Range("H86").Value = 6199
$ws.Range("I86").Value = 2038
$ws.Range("K86").Value = 2038
$ws.Range("M86").Value = -915
$ws.Range("H89").Value = 6199
$ws.Range("I89").Value = 2038
$ws.Range("K89").Value = 10190
$ws.Range("M89").Value = -4574
$ws.Range("H94").Value = 1980.5
$ws.Range("I94").Value = 1980.5
$ws.Range("K94").Value = 1980.5
$ws.Range("M94").Value = -1529.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2151.5715
$ws.Range("I16").Value = 2151.5715
$ws.Range("K16").Value = 2151.5715
$ws.Range("M16").Value = -1864.5715
$ws.Range("H22").Value = 2070.1
$ws.Range("J22").Value = 3056.6
$ws.Range("L22").Value = 3056.6
$ws.Range("N22").Value = -3756.6
$ws.Range("H28").Value = 20627.334
$ws.Range("J28").Value = 20627.334
$ws.Range("L28").Value = 20627.334
$ws.Range("N28").Value = -21117.334
$ws.Range("H99").Value = 2484
$ws.Range("I99").Value = 2482.5
$ws.Range("J99").Value = 2485
$ws.Range("K99").Value = 2482.5
$ws.Range("L99").Value = 2485
$ws.Range("M99").Value = -984.5
$ws.Range("N99").Value = -5481
$ws.Range("H113").Value = 2151.5715
$ws.Range("I113").Value = 2151.5715
$ws.Range("K113").Value = 2151.5715
$ws.Range("M113").Value = 18.42849999999999
$ws.Range("H126").Value = 2484
$ws.Range("I126").Value = 2482.5
$ws.Range("J126").Value = 2485
$ws.Range("K126").Value = 7447.5
$ws.Range("L126").Value = 7455
$ws.Range("M126").Value = -4977.5
$ws.Range("N126").Value = -12395

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 190
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H51").Value = 500
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H59").Value = 825
$ws.Range("H80").Value = 4279.5835
$ws.Range("I80").Value = 4020.4736
$ws.Range("J80").Value = 5264.2
$ws.Range("K80").Value = 12061.4208
$ws.Range("L80").Value = 15792.6
$ws.Range("M80").Value = -11125.4208
$ws.Range("N80").Value = -17664.6
$ws.Range("H83").Value = 4279.5835
$ws.Range("I83").Value = 4020.4736
$ws.Range("J83").Value = 5264.2
$ws.Range("K83").Value = 36184.2624
$ws.Range("L83").Value = 47377.8
$ws.Range("M83").Value = -31504.2624
$ws.Range("N83").Value = -56737.8
$ws.Range("H134").Value = 1916.6666
$ws.Range("I134").Value = 1916.6666
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5749.9998
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -679.9997999999996
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 29944
$ws.Range("J39").Value = 29944
$ws.Range("L39").Value = 29944
$ws.Range("N39").Value = -31008
$ws.Range("H80").Value = 3997
$ws.Range("I80").Value = 4329.6665
$ws.Range("J80").Value = 2999
$ws.Range("K80").Value = 4329.6665
$ws.Range("L80").Value = 2999
$ws.Range("M80").Value = -3331.6665
$ws.Range("N80").Value = -4995
$ws.Range("H83").Value = 3997
$ws.Range("I83").Value = 4329.6665
$ws.Range("J83").Value = 2999
$ws.Range("K83").Value = 21648.3325
$ws.Range("L83").Value = 14995
$ws.Range("M83").Value = -16656.3325
$ws.Range("N83").Value = -24979
$ws.Range("H98").Value = 11761.8
$ws.Range("J98").Value = 11761.8
$ws.Range("L98").Value = 11761.8
$ws.Range("N98").Value = -17751.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2383.158
$ws.Range("I40").Value = 3040
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 3040
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -2904
$ws.Range("N40").Value = -2272
$ws.Range("H61").Value = 166671730
$ws.Range("I61").Value = 200005070
$ws.Range("K61").Value = 200005070
$ws.Range("M61").Value = -200004868
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H93").Value = 2632
$ws.Range("I93").Value = 2715.1428
$ws.Range("K93").Value = 2715.1428
$ws.Range("M93").Value = -1467.1428
$ws.Range("H113").Value = 166671730
$ws.Range("I113").Value = 200005070
$ws.Range("K113").Value = 200005070
$ws.Range("M113").Value = -200002900
$ws.Range("H140").Value = 27124.75
$ws.Range("J140").Value = 27124.75
$ws.Range("L140").Value = 27124.75
$ws.Range("N140").Value = -37484.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 16813
$ws.Range("I34").Value = 16813
$ws.Range("K34").Value = 16813
$ws.Range("M34").Value = -16610
$ws.Range("H112").Value = 20462.334
$ws.Range("J112").Value = 20462.334
$ws.Range("L112").Value = 20462.334
$ws.Range("N112").Value = -23416.334
$ws.Range("H136").Value = 2213.8518
$ws.Range("I136").Value = 1372.05
$ws.Range("K136").Value = 4116.15
$ws.Range("M136").Value = -1566.15
